$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (prices, 1h volume %, and a few re-ranked rows)

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '62.784.03'
$ws.Cells.Item(2, 5).Value = '  -6.34%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.252.17'
$ws.Cells.Item(3, 5).Value = '  -7.45%  '

$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.28%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '178.70'
$ws.Cells.Item(5, 5).Value = '  -11.98%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '514.58'
$ws.Cells.Item(6, 5).Value = '  -6.68%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.595'
$ws.Cells.Item(7, 5).Value = '  -0.86%  '

$ws.Cells.Item(8, 2).Value = 'LidoStakedEther'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '3.250.55'
$ws.Cells.Item(8, 5).Value = '  -7.32%  '

$ws.Cells.Item(9, 2).Value = 'USDC'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '1.00'
$ws.Cells.Item(9, 5).Value = '  +0.01%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.614'
$ws.Cells.Item(10, 5).Value = '  -6.17%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '57.50'
$ws.Cells.Item(11, 5).Value = '  -6.02%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.130'
$ws.Cells.Item(12, 5).Value = '  -8.91%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.0000253'
$ws.Cells.Item(13, 5).Value = '  -6.45%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '9.08'
$ws.Cells.Item(14, 5).Value = '  -7.89%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '3.753.67'
$ws.Cells.Item(15, 5).Value = '  -8.02%  '

$ws.Cells.Item(16, 5).Value = '  -6.54%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '3.244.11'
$ws.Cells.Item(17, 5).Value = '  -7.84%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '17.54'
$ws.Cells.Item(18, 5).Value = '  -5.58%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '62.591.74'
$ws.Cells.Item(19, 5).Value = '  -6.45%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '10.89'
$ws.Cells.Item(20, 5).Value = '  -7.99%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '0.944'
$ws.Cells.Item(21, 5).Value = '  -9.05%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '369.32'
$ws.Cells.Item(22, 5).Value = '  -5.11%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '11.21'
$ws.Cells.Item(23, 5).Value = '  -5.86%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '3.68'
$ws.Cells.Item(24, 5).Value = '  -8.30%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '79.28'
$ws.Cells.Item(25, 5).Value = '  -4.02%  '

$ws.Cells.Item(26, 2).Value = 'LEO'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '6.05'
$ws.Cells.Item(26, 5).Value = '  -1.37%  '

$ws.Cells.Item(27, 2).Value = 'Toncoin'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '3.79'
$ws.Cells.Item(27, 5).Value = '  +1.70%  '

$ws.Cells.Item(28, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '11.32'
$ws.Cells.Item(28, 5).Value = '  -5.67%  '

$ws.Cells.Item(29, 2).Value = 'ImmutableX'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.60'
$ws.Cells.Item(29, 5).Value = '  -7.12%  '

$ws.Cells.Item(30, 2).Value = 'Filecoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '8.29'
$ws.Cells.Item(30, 5).Value = '  -6.59%  '

$ws.Cells.Item(31, 2).Value = 'EthereumClassic'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '28.32'
$ws.Cells.Item(31, 5).Value = '  -7.68%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '6.70'
$ws.Cells.Item(32, 5).Value = '  -8.24%  '

$ws.Cells.Item(33, 2).Value = 'Bittensor'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '630.69'
$ws.Cells.Item(33, 5).Value = '  -8.87%  '

$ws.Cells.Item(34, 2).Value = 'Cosmos'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '11.15'
$ws.Cells.Item(34, 5).Value = '  -5.02%  '

$ws.Cells.Item(35, 2).Value = 'Hedera'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.105'
$ws.Cells.Item(35, 5).Value = '  -4.86%  '

$ws.Cells.Item(36, 2).Value = 'OKB'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '58.23'
$ws.Cells.Item(36, 5).Value = '  -7.69%  '

$ws.Cells.Item(37, 2).Value = 'Dai'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.00'
$ws.Cells.Item(37, 5).Value = '  -0.02%  '

$ws.Cells.Item(38, 2).Value = 'TheGraph'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.393'
$ws.Cells.Item(38, 5).Value = '  -3.17%  '

$ws.Cells.Item(39, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '36.11'
$ws.Cells.Item(39, 5).Value = '  -9.56%  '

$ws.Cells.Item(40, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.997'
$ws.Cells.Item(40, 5).Value = '  -0.26%  '

$ws.Cells.Item(41, 2).Value = 'Maker'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '2.937.38'
$ws.Cells.Item(41, 5).Value = '  -5.74%  '

$ws.Cells.Item(42, 2).Value = 'Kaspa'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.123'
$ws.Cells.Item(42, 5).Value = '  -4.58%  '

$ws.Cells.Item(43, 2).Value = 'PEPE'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.0₃0650'
$ws.Cells.Item(43, 5).Value = '  -8.19%  '

$ws.Cells.Item(44, 2).Value = 'Fetch.AI'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '2.43'
$ws.Cells.Item(44, 5).Value = '  -3.30%  '

$ws.Cells.Item(45, 2).Value = 'ThetaToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '2.65'
$ws.Cells.Item(45, 5).Value = '  -13.94%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.0389'
$ws.Cells.Item(46, 5).Value = '  -2.58%  '

$ws.Cells.Item(47, 2).Value = 'WEMIXToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.58'
$ws.Cells.Item(47, 5).Value = '  -5.12%  '

$ws.Cells.Item(48, 2).Value = 'Stacks'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '2.79'
$ws.Cells.Item(48, 5).Value = '  +6.10%  '

$ws.Cells.Item(49, 2).Value = 'Stellar'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.124'
$ws.Cells.Item(49, 5).Value = '  -2.53%  '

$ws.Cells.Item(50, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '2.91'
$ws.Cells.Item(50, 5).Value = '  -1.37%  '

$ws.Cells.Item(51, 2).Value = 'dogwifhat'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '2.45'
$ws.Cells.Item(51, 5).Value = '  -14.21%  '
